$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 267
$ws.Range("J2").Value = 402
$ws.Range("L2").Value = 402
$ws.Range("N2").Value = -628
$ws.Range("H19").Value = 803.6316
$ws.Range("I19").Value = 349
$ws.Range("J19").Value = 1134.2727
$ws.Range("K19").Value = 349
$ws.Range("L19").Value = 1134.2727
$ws.Range("M19").Value = -174
$ws.Range("N19").Value = -1484.2727
$ws.Range("H32").Value = 438.46155
$ws.Range("J32").Value = 505.22223
$ws.Range("L32").Value = 505.22223
$ws.Range("N32").Value = -1157.22223
$ws.Range("H51").Value = 7208.375
$ws.Range("I51").Value = 14211.223
$ws.Range("J51").Value = 3006.6667
$ws.Range("K51").Value = 14211.223
$ws.Range("L51").Value = 3006.6667
$ws.Range("M51").Value = -13727.223
$ws.Range("N51").Value = -3974.6667
$ws.Range("H64").Value = 58511.11
$ws.Range("I64").Value = 93618.17999999999
$ws.Range("J64").Value = 3342.8572
$ws.Range("K64").Value = 93618.17999999999
$ws.Range("L64").Value = 3342.8572
$ws.Range("M64").Value = -93370.17999999999
$ws.Range("N64").Value = -3838.8572
$ws.Range("H67").Value = 58511.11
$ws.Range("I67").Value = 93618.17999999999
$ws.Range("J67").Value = 3342.8572
$ws.Range("K67").Value = 93618.17999999999
$ws.Range("L67").Value = 3342.8572
$ws.Range("M67").Value = -92760.17999999999
$ws.Range("N67").Value = -5058.8572
$ws.Range("H107").Value = 1097.8334
$ws.Range("I107").Value = 1097.8334
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1097.8334
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 822.1666
$ws.Range("N107").ClearContents()
$ws.Range("H132").Value = 3627650.8
$ws.Range("I132").Value = 4634630
$ws.Range("J132").Value = 2525.3333
$ws.Range("K132").Value = 13903890
$ws.Range("L132").Value = 7575.999899999999
$ws.Range("M132").Value = -13901360
$ws.Range("N132").Value = -12635.9999
$ws.Range("H135").Value = 1280.5193
$ws.Range("I135").Value = 499.8889
$ws.Range("J135").Value = 3036.9375
$ws.Range("K135").Value = 4499.0001
$ws.Range("L135").Value = 27332.4375
$ws.Range("M135").Value = -1964.0001
$ws.Range("N135").Value = -32402.4375
$ws.Range("H138").Value = 2286.1
$ws.Range("I138").Value = 1509.8636
$ws.Range("J138").Value = 2505.0386
$ws.Range("K138").Value = 4529.5908
$ws.Range("L138").Value = 7515.1158
$ws.Range("M138").Value = 610.4092000000001
$ws.Range("N138").Value = -17795.1158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 33433.055
$ws.Range("I32").Value = 14845.917
$ws.Range("J32").Value = 69120.36
$ws.Range("K32").Value = 14845.917
$ws.Range("L32").Value = 69120.36
$ws.Range("M32").Value = -14558.917
$ws.Range("N32").Value = -69694.36
$ws.Range("H74").Value = 949.4375
$ws.Range("I74").Value = 636.0909
$ws.Range("K74").Value = 636.0909
$ws.Range("M74").Value = 237.9091
$ws.Range("H77").Value = 949.4375
$ws.Range("I77").Value = 636.0909
$ws.Range("K77").Value = 3180.4545
$ws.Range("M77").Value = 1187.5455
$ws.Range("H97").Value = 41199.24
$ws.Range("I97").Value = 48596.57
$ws.Range("K97").Value = 48596.57
$ws.Range("M97").Value = -48100.57

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 21920.4
$ws.Range("J35").Value = 21920.4
$ws.Range("L35").Value = 21920.4
$ws.Range("N35").Value = -22540.4
$ws.Range("H82").Value = 20330.715
$ws.Range("I82").Value = 3164
$ws.Range("J82").Value = 33205.75
$ws.Range("K82").Value = 3164
$ws.Range("L82").Value = 33205.75
$ws.Range("M82").Value = -2781
$ws.Range("N82").Value = -33971.75
$ws.Range("H85").Value = 20330.715
$ws.Range("I85").Value = 3164
$ws.Range("J85").Value = 33205.75
$ws.Range("K85").Value = 3164
$ws.Range("L85").Value = 33205.75
$ws.Range("M85").Value = -1838
$ws.Range("N85").Value = -35857.75
$ws.Range("H134").Value = 3739.025
$ws.Range("I134").Value = 4045.2222
$ws.Range("J134").Value = 3103.077
$ws.Range("K134").Value = 12135.6666
$ws.Range("L134").Value = 9309.231
$ws.Range("M134").Value = -9600.6666
$ws.Range("N134").Value = -14379.231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 221.75
$ws.Range("I7").Value = 145.66667
$ws.Range("K7").Value = 145.66667
$ws.Range("M7").Value = -32.66667000000001
$ws.Range("H10").Value = 500
$ws.Range("I10").Value = 500
$ws.Range("K10").Value = 500
$ws.Range("M10").Value = -361
$ws.Range("H31").Value = 40266.41
$ws.Range("I31").Value = 1956.4445
$ws.Range("J31").Value = 51759.4
$ws.Range("K31").Value = 1956.4445
$ws.Range("L31").Value = 51759.4
$ws.Range("M31").Value = -1661.4445
$ws.Range("N31").Value = -52349.4
$ws.Range("H34").Value = 40266.41
$ws.Range("I34").Value = 1956.4445
$ws.Range("J34").Value = 51759.4
$ws.Range("K34").Value = 1956.4445
$ws.Range("L34").Value = 51759.4
$ws.Range("M34").Value = -1754.4445
$ws.Range("N34").Value = -52163.4
$ws.Range("H58").Value = 1638.1951
$ws.Range("I58").Value = 1507.6
$ws.Range("J58").Value = 1842.25
$ws.Range("K58").Value = 1507.6
$ws.Range("L58").Value = 1842.25
$ws.Range("M58").Value = -1304.6
$ws.Range("N58").Value = -2248.25
$ws.Range("H134").Value = 1570.4736
$ws.Range("I134").Value = 840.1
$ws.Range("J134").Value = 2382
$ws.Range("K134").Value = 2520.3
$ws.Range("L134").Value = 7146
$ws.Range("M134").Value = 14.69999999999982
$ws.Range("N134").Value = -12216
$ws.Range("H136").Value = 1638.1951
$ws.Range("I136").Value = 1507.6
$ws.Range("J136").Value = 1842.25
$ws.Range("K136").Value = 4522.799999999999
$ws.Range("L136").Value = 5526.75
$ws.Range("M136").Value = -1972.799999999999
$ws.Range("N136").Value = -10626.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1334
$ws.Range("I3").Value = 1334
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 4002
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -3890
$ws.Range("N3").ClearContents()
$ws.Range("H12").Value = 54.833332
$ws.Range("J12").Value = 58.8125
$ws.Range("L12").Value = 176.4375
$ws.Range("N12").Value = -522.4375
$ws.Range("H14").Value = 740.5
$ws.Range("I14").Value = 740.5
$ws.Range("K14").Value = 2221.5
$ws.Range("M14").Value = -2048.5
$ws.Range("H33").Value = 1435.4584
$ws.Range("I33").Value = 1230.3636
$ws.Range("J33").Value = 1609
$ws.Range("K33").Value = 7382.1816
$ws.Range("L33").Value = 9654
$ws.Range("M33").Value = -7099.1816
$ws.Range("N33").Value = -10220
$ws.Range("H86").Value = 1149.8334
$ws.Range("I86").Value = 1500
$ws.Range("J86").Value = 974.75
$ws.Range("K86").Value = 4500
$ws.Range("L86").Value = 2924.25
$ws.Range("M86").Value = -3314
$ws.Range("N86").Value = -5296.25
$ws.Range("H89").Value = 1149.8334
$ws.Range("I89").Value = 1500
$ws.Range("J89").Value = 974.75
$ws.Range("K89").Value = 13500
$ws.Range("L89").Value = 8772.75
$ws.Range("M89").Value = -7572
$ws.Range("N89").Value = -20628.75
$ws.Range("H131").Value = 634682.0600000001
$ws.Range("J131").Value = 716194.0600000001
$ws.Range("L131").Value = 2148582.18
$ws.Range("N131").Value = -2158662.18
$ws.Range("H132").Value = 1793.2069
$ws.Range("J132").Value = 2120.45
$ws.Range("L132").Value = 19084.05
$ws.Range("N132").Value = -24144.05

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 12499.923
$ws.Range("J46").Value = 12499.923
$ws.Range("L46").Value = 12499.923
$ws.Range("N46").Value = -12811.923

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H80").Value = 10259.72
$ws.Range("J80").Value = 10259.72
$ws.Range("L80").Value = 10259.72
$ws.Range("N80").Value = -12505.72
$ws.Range("H83").Value = 10259.72
$ws.Range("J83").Value = 10259.72
$ws.Range("L83").Value = 30779.16
$ws.Range("N83").Value = -42011.16
$ws.Range("H132").Value = 5105.3
$ws.Range("I132").Value = 7675.3335
$ws.Range("J132").Value = 3002.5454
$ws.Range("K132").Value = 23026.0005
$ws.Range("L132").Value = 9007.636200000001
$ws.Range("M132").Value = -20496.0005
$ws.Range("N132").Value = -14067.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3348.7542
$ws.Range("I132").Value = 1657.44
$ws.Range("J132").Value = 11036.546
$ws.Range("K132").Value = 4972.32
$ws.Range("L132").Value = 33109.638
$ws.Range("M132").Value = -2442.32
$ws.Range("N132").Value = -38169.638
$ws.Range("H136").Value = 16971.078
$ws.Range("I136").Value = 28284.416
$ws.Range("J136").Value = 2425.3572
$ws.Range("K136").Value = 84853.24800000001
$ws.Range("L136").Value = 7276.071599999999
$ws.Range("M136").Value = -82303.24800000001
$ws.Range("N136").Value = -12376.0716
